# Edit slide 42 ("Abstract Syntax Trees (continued)") of the presentation:
#  - Rework Example 2's intro sentence and grammar rule (while/loop terminals)
#  - Rework the explanatory paragraph to reference the new "while"/"loop" terminals
#  - Nudge the small diagram (Group 1) down slightly

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(42)

# --- Body placeholder (the shape with the worked example) ---
$bodyShape = $null
foreach ($sh in $s.Shapes) {
    if ($sh.Name -eq "Rectangle 3") {
        $bodyShape = $sh
    }
}

$tf = $bodyShape.TextFrame
$tr = $tf.TextRange

# 1) First sentence of paragraph 1: split off "while" into Consolas font.
$run1 = $tr.Paragraphs(1).Runs(1)
$newIntro = "Example 2: Consider the following rule for a while statement."
$run1.Text = $newIntro
$introStart = $run1.Start
$idxWhile = $newIntro.IndexOf("while")
$introWhile = $tf.TextRange.Characters($introStart + $idxWhile, 5)
$introWhile.Font.Name = "Consolas"

# 2) Grammar rule line: " = "while" "(" " -> " = "while" "
$found = $tf.TextRange.Find(' = "while" "(" ')
$found.Text = ' = "while" '

# 3) Grammar rule line: ' ")" statement .' -> ' "loop" statement .'
$found2 = $tf.TextRange.Find(' ")" statement .')
$found2.Text = ' "loop" statement .'

# 4) Paragraph 2: mention the "while" and "loop" terminal symbols explicitly.
$para2run1 = $tf.TextRange.Paragraphs(2).Runs(1)
$apos = [char]0x2019
$newPara2 = "Once a while statement has been parsed, we don" + $apos + "t need to retain the terminal symbols ""while"" and ""loop"".  The abstract syntax tree for a while statement would contain only "
$para2run1.Text = $newPara2
$para2Start = $para2run1.Start

$idxWhile2 = $newPara2.IndexOf([string]([char]34) + "while" + [char]34)
$whileRange = $tf.TextRange.Characters($para2Start + $idxWhile2, 7)
$whileRange.Font.Name = "Consolas"

$idxLoop2 = $newPara2.IndexOf([string]([char]34) + "loop" + [char]34)
$loopRange = $tf.TextRange.Characters($para2Start + $idxLoop2, 6)
$loopRange.Font.Name = "Consolas"

# --- Move the small "Group 1" diagram down a bit ---
$groupShape = $null
foreach ($sh in $s.Shapes) {
    if ($sh.Name -eq "Group 1") {
        $groupShape = $sh
    }
}
$groupShape.Top = 346.44472440944884
